$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------------
# Source data (before the edit):
#   Row 243 = HYN / Taizhou, China   (Asia / Taizhou / China / CN / blank lat / blank lon)
#   Row 260 = XNN / Xining, China    (row to be removed entirely)
#   Row 272 = IAD / Ashburn, VA, United States
#
# Target data (after the edit):
#   - The XNN (Xining, China) row is deleted outright.
#   - The HYN (Taizhou, China) row is relocated so that it sits immediately
#     before the IAD (Ashburn) row instead of its original spot.
#   - Dimension shrinks from A1:H331 to A1:H330 (net loss of exactly one row).
# ----------------------------------------------------------------------------

# 1) Delete the XNN (Xining, China) row completely; everything below shifts up by one.
$ws.Rows("260:260").Delete()

# After the delete above, HYN is still at row 243, and IAD (previously row 272)
# is now at row 271.

# 2) Remove the HYN row from its original location; everything below shifts up again.
$ws.Rows("243:243").Delete()

# IAD is now at row 270.

# 3) Insert a fresh blank row right above the new IAD position (row 270), which
#    is where HYN needs to end up.
$ws.Rows("270:270").Insert()

# 4) Pick up the formatting (styles/borders/fonts) of the neighboring row so the
#    newly inserted row matches the sheet's existing look (bold/boxed colo column).
$ws.Range("A269:H269").Copy()
$ws.Range("A270:H270").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 5) Populate the relocated HYN row with its original values. Latitude/longitude
#    stay blank, matching the source data.
$ws.Range("A270").Value = "HYN"
$ws.Range("B270").Value = "Taizhou, China"
$ws.Range("C270").Value = "Asia"
$ws.Range("D270").Value = "Taizhou"
$ws.Range("E270").Value = "China"
$ws.Range("F270").Value = "CN"
